# "Generate Report for Handback" - refresh the handback timestamps that were
# regenerated for the 5fb35e66-... (in-sync) report row (row 3 in each table).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G3 "Latest HO Xliff Generate Date" (same value also shown on de-de!H3)
$wsOverview.Range("G3").Value = "2016-09-08 05:00:22"

# zh-cn table, row 3: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-09-08 05:00:00"
$wsZhCn.Range("K3").Value = "2016-09-08 05:00:42"

# de-de table, row 3: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-09-08 05:00:22"
$wsDeDe.Range("K3").Value = "2016-09-08 05:00:51"
